$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Demographic")
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 33
